$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Tanuki row's sprite name was renamed from "tanuki_mario" to "tanuki".
# Do this first so the new shared string ends up in the same relative slot
# the author's workbook has it in (right after the other renamed strings).
$ws.Range("B3").Value = "tanuki"

# The sheet previously ended with a lone "END" marker on row 5 (column A).
# Two new enemy rows (Ladon, Lilim) are being inserted before that marker,
# so first push the existing row 5 content down to row 7.
$ws.Cells.Item(7, 1).Value = $ws.Cells.Item(5, 1).Value2

# Row 5: Ladon - stat line copied from the Tanuki row (row 3), just a new
# Name(ID) / sprite_name pair.
$ws.Cells.Item(5, 1).Value  = "Ladon"
$ws.Cells.Item(5, 2).Value  = "ladon"
$ws.Cells.Item(5, 3).Value  = 100
$ws.Cells.Item(5, 4).Value  = 0
$ws.Cells.Item(5, 5).Value  = 3
$ws.Cells.Item(5, 6).Value  = 0.95
$ws.Cells.Item(5, 7).Value  = -0.2
$ws.Cells.Item(5, 8).Value  = 2.5
$ws.Cells.Item(5, 9).Value  = 1
$ws.Cells.Item(5, 10).Value = 10
$ws.Cells.Item(5, 11).Value = "Attacker1"
$ws.Cells.Item(5, 12).Value = "none"
$ws.Cells.Item(5, 13).Value = 1.2
$ws.Cells.Item(5, 14).Value = 1
$ws.Cells.Item(5, 15).Value = 1
$ws.Cells.Item(5, 16).Value = 1
$ws.Cells.Item(5, 17).Value = "GROUP/DEFAULT"
$ws.Cells.Item(5, 18).Value = "spear"
$ws.Cells.Item(5, 19).Value = "null"
$ws.Cells.Item(5, 20).Value = "null"
$ws.Cells.Item(5, 21).Value = "spear"
$ws.Cells.Item(5, 22).Value = "null"
$ws.Cells.Item(5, 23).Value = "aimed"
$ws.Cells.Item(5, 24).Value = "END"

# Row 6: Lilim - stat line copied from "The Evil Eye" row (row 4), just a
# new Name(ID) / sprite_name pair.
$ws.Cells.Item(6, 1).Value  = "Lilim"
$ws.Cells.Item(6, 2).Value  = "lilim"
$ws.Cells.Item(6, 3).Value  = 50
$ws.Cells.Item(6, 4).Value  = 0
$ws.Cells.Item(6, 5).Value  = 2
$ws.Cells.Item(6, 6).Value  = 1.25
$ws.Cells.Item(6, 7).Value  = 0.2
$ws.Cells.Item(6, 8).Value  = 1.5
$ws.Cells.Item(6, 9).Value  = 1
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = "HealthLow1"
$ws.Cells.Item(6, 12).Value = "75,25"
$ws.Cells.Item(6, 13).Value = -1
$ws.Cells.Item(6, 14).Value = 2
$ws.Cells.Item(6, 15).Value = 0.75
$ws.Cells.Item(6, 16).Value = -2
$ws.Cells.Item(6, 17).Value = "GROUP/DEFAULT"
$ws.Cells.Item(6, 18).Value = "sword"
$ws.Cells.Item(6, 19).Value = "null"
$ws.Cells.Item(6, 20).Value = "null"
$ws.Cells.Item(6, 21).Value = "spear"
$ws.Cells.Item(6, 22).Value = "fire"
$ws.Cells.Item(6, 23).Value = "null"
$ws.Cells.Item(6, 24).Value = "GROUP/HEALTH_LOW"
$ws.Cells.Item(6, 25).Value = "quake"
$ws.Cells.Item(6, 26).Value = "null"
$ws.Cells.Item(6, 27).Value = "null"
$ws.Cells.Item(6, 28).Value = "END"

# Selection moved to C6 (wherever the author's cursor ended up).
[void]$ws.Range("C6").Select()
